# Updated cryptos list on Fri Jun  2 18:28:39 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.068.63'
$ws.Range("E2").Value = '  +0.22%  '
$ws.Range("D3").Value = '1.893.96'
$ws.Range("E3").Value = '  +1.03%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.09'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.49%  '
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5187'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3770'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.72%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07223'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.14'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8919'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.49%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07664'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.85%  '
$ws.Range("D13").Value = '1.901.54'
$ws.Range("E13").Value = '  +3.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.18'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.97%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.225'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.56%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008509'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.50'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.58%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.13%  '
$ws.Range("D20").Value = '27.124.26'
$ws.Range("E20").Value = '  +0.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.062'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.61%  '
$ws.Range("D22").Value = '2.151.31'
$ws.Range("E22").Value = '  +2.96%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.57'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.61%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.291'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +9.56%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '146.19'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.48%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.733'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.17%  '
$ws.Range("E28").Value = '  +0.78%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.45'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.88%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.963'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.786'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.32%  '
$ws.Range("E32").Value = '  +0.49%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05047'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.51%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.235'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.32%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7757'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.975'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.292'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.82%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.582'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.45%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5614'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.43%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01987'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.072'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.27%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.983'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.61%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.627'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.09%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '118.67'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.69%  '
$ws.Range("E45").Value = '  +2.58%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4819'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.62%  '
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.001'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.18%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.15'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.35%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.593'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '37.50'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.60%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.94'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.06%  '
